$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.946.24"
$ws.Range("E2").Value = "  +0.70%  "

$ws.Range("D3").Value = "'1.593.92"
$ws.Range("E3").Value = "  +0.77%  "

$ws.Range("E4").Value = "  +0.34%  "

$ws.Range("D5").Value = "'210.32"
$ws.Range("E5").Value = "  +0.60%  "

$ws.Range("E6").Value = "  +0.37%  "

$ws.Range("E7").Value = "  -0.15%  "

$ws.Range("E9").Value = "  -1.37%  "

$ws.Range("D10").Value = "'18.00"
$ws.Range("E10").Value = "  -1.14%  "

$ws.Range("E11").Value = "  +2.91%  "

$ws.Range("D12").Value = "'1.818.14"
$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("D13").Value = "'1.629.08"
$ws.Range("E13").Value = "  +2.99%  "

$ws.Range("D14").Value = "'3.99"
$ws.Range("E14").Value = "  -0.92%  "

$ws.Range("D15").Value = "'0.511"
$ws.Range("E15").Value = "  -0.45%  "

$ws.Range("D16").Value = "'25.973.44"
$ws.Range("E16").Value = "  +0.82%  "

$ws.Range("D17").Value = "'59.97"
$ws.Range("E17").Value = "  -0.44%  "

$ws.Range("D18").Value = "'0.0₃0720"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("E19").Value = "  +0.37%  "

$ws.Range("D20").Value = "'199.71"
$ws.Range("E20").Value = "  +4.20%  "

$ws.Range("E21").Value = "  +0.54%  "

$ws.Range("D22").Value = "'9.22"
$ws.Range("E22").Value = "  -2.18%  "

$ws.Range("D23").Value = "'5.99"
$ws.Range("E23").Value = "  +0.88%  "

$ws.Range("D24").Value = "'1.79"
$ws.Range("E24").Value = "  +4.90%  "

$ws.Range("D25").Value = "'141.73"
$ws.Range("E25").Value = "  +0.54%  "

$ws.Range("E26").Value = "  +0.47%  "

$ws.Range("E27").Value = "  -8.33%  "

$ws.Range("D28").Value = "'15.06"
$ws.Range("E28").Value = "  -0.79%  "

$ws.Range("D29").Value = "'6.44"
$ws.Range("E29").Value = "  -0.48%  "

$ws.Range("E30").Value = "  +0.38%  "

$ws.Range("D31").Value = "'0.0474"
$ws.Range("E31").Value = "  +0.32%  "

$ws.Range("D32").Value = "'3.09"
$ws.Range("E32").Value = "  -0.95%  "

$ws.Range("E33").Value = "  -2.69%  "

$ws.Range("E34").Value = "  -1.94%  "

$ws.Range("E35").Value = "  +2.29%  "

$ws.Range("D36").Value = "'1.122.86"
$ws.Range("E36").Value = "  +1.89%  "

$ws.Range("D37").Value = "'0.0161"
$ws.Range("E37").Value = "  +6.77%  "

$ws.Range("E38").Value = "  +0.37%  "

$ws.Range("E39").Value = "  -0.52%  "

$ws.Range("D40").Value = "'0.780"
$ws.Range("E40").Value = "  -0.84%  "

$ws.Range("E41").Value = "  -3.16%  "

$ws.Range("D42").Value = "'0.778"
$ws.Range("E42").Value = "  -3.62%  "

$ws.Range("D43").Value = "'1.729.21"
$ws.Range("E43").Value = "  +0.73%  "

$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'92.45"
$ws.Range("E44").Value = "  -1.13%  "

$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").Value = "'5.08"
$ws.Range("E45").Value = "  -1.35%  "

$ws.Range("D46").Value = "'1.49"
$ws.Range("E46").Value = "  -0.50%  "

$ws.Range("D47").Value = "'53.16"
$ws.Range("E47").Value = "  -0.22%  "

$ws.Range("E48").Value = "  -1.39%  "

$ws.Range("E49").Value = "  +0.12%  "

$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("D51").Value = "'0.0₇0914"
$ws.Range("E51").Value = "  -17.91%  "
